$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116, shifting existing rows 116-136 down to 117-137
$ws.Rows("116").Insert()

# Populate the newly inserted row 116 with the new weekly price observation
$ws.Range("A116").Value = 6
$ws.Range("B116").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C116").Value = "Metropolitana"
$ws.Range("D116").Value = 44474
$ws.Range("E116").Value = 13
$ws.Range("F116").Value = 100112022
$ws.Range("G116").Value = "Arveja Verde"
$ws.Range("H116").Value = "Sin especificar"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 180
$ws.Range("K116").Value = 23000
$ws.Range("L116").Value = 24000
$ws.Range("M116").Value = 23444
$ws.Range("N116").Value = "$/malla 25 kilos"
$ws.Range("O116").Value = "Provincia de Huasco"
$ws.Range("P116").Value = 938
$ws.Range("Q116").Value = 25
$ws.Range("R116").Value = "Hortaliza"
